$wb = $excel.ActiveWorkbook

$tocSheet = $wb.Worksheets.Item("!!_Table of contents")
$schemaSheet = $wb.Worksheets.Item("!!_Schema")
$txSheet = $wb.Worksheets.Item("!!Transaction")

# Sheets are protected (no password) - unprotect before editing
$tocSheet.Unprotect()
$schemaSheet.Unprotect()
$txSheet.Unprotect()

# Update ObjTables version/date headers on the Table of contents sheet
$tocSheet.Range("A1").Value = "!!!ObjTables objTablesVersion='1.0.0' date='2020-05-29 00:18:57'"
$tocSheet.Range("A2").Value = "!!ObjTables type='TableOfContents' tableFormat='row' description='Table of contents' date='2020-05-29 00:18:57' objTablesVersion='1.0.0'"

# Update the Schema sheet header and add Verbose name values for attributes
$schemaSheet.Range("A1").Value = "!!ObjTables type='Schema' tableFormat='row' description='Table/model and column/attribute definitions' date='2020-05-29 00:18:57' objTablesVersion='1.0.0'"

$schemaSheet.Range("E4").Value = "Amount"
$schemaSheet.Range("E5").Value = "Category"
$schemaSheet.Range("E6").Value = "Date"
$schemaSheet.Range("E7").Value = "Payee"

# Update the Transaction data sheet header
$txSheet.Range("A1").Value = "!!ObjTables type='Data' tableFormat='row' class='Transaction' name='Transaction' description='Stores transactions' date='2020-05-29 00:18:57' objTablesVersion='1.0.0'"

# Restore sheet protection (matches original: sheet/objects/scenarios protected;
# inserting & deleting rows explicitly allowed; everything else left at the
# protected default).
$tocSheet.Protect($null, $true, $true, $true, $false, $false, $false, $false, $false, $true, $false, $false, $true, $false, $false, $false)
$schemaSheet.Protect($null, $true, $true, $true, $false, $false, $false, $false, $false, $true, $false, $false, $true, $false, $false, $false)
$txSheet.Protect($null, $true, $true, $true, $false, $false, $false, $false, $false, $true, $false, $false, $true, $false, $false, $false)
